$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.582.66"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.352.74"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'557.68"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").Value = "'131.71"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -3.13%  "
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "'5.58"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("E13").Value = "  -5.68%  "
$ws.Range("D14").Value = "2.769.81"
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "59.559.42"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "2.351.21"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "'4.42"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'318.88"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").Value = "'6.59"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'63.94"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'8.31"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "'171.47"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "0.0₃0743"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  +6.75%  "
$ws.Range("D33").Value = "'0.398"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").Value = "'312.90"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").Value = "'144.26"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "'3.46"
$ws.Range("E43").Value = "  -5.45%  "
$ws.Range("D44").Value = "'0.0955"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "'18.79"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").Value = "'0.0215"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "'11.05"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -0.36%  "
